# Auto-generated Excel COM-interop script
# Applies numeric updates to the Leve profit-tracking tables across all class sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled price-refresh run.

$wb = $excel.ActiveWorkbook

# ALC row 75
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 113500
$ws.Range("J75").Value = 113500
$ws.Range("L75").Value = 113500
$ws.Range("N75").Value = -115372

# ALC row 78
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H78").Value = 113500
$ws.Range("J78").Value = 113500
$ws.Range("L78").Value = 340500
$ws.Range("N78").Value = -349860

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4927.5557
$ws.Range("I86").Value = 6153.3335
$ws.Range("J86").Value = 4314.6665
$ws.Range("K86").Value = 6153.3335
$ws.Range("L86").Value = 4314.6665
$ws.Range("M86").Value = -5030.3335
$ws.Range("N86").Value = -6560.6665

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4927.5557
$ws.Range("I89").Value = 6153.3335
$ws.Range("J89").Value = 4314.6665
$ws.Range("K89").Value = 30766.6675
$ws.Range("L89").Value = 21573.3325
$ws.Range("M89").Value = -25150.6675
$ws.Range("N89").Value = -32805.3325

# ALC row 101
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 3627.8572
$ws.Range("J101").Value = 3889
$ws.Range("L101").Value = 11667
$ws.Range("N101").Value = -14911

# ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 5000
$ws.Range("J111").Value = 5000
$ws.Range("L111").Value = 15000
$ws.Range("N111").Value = -21134

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2004.5454

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 5428.9756
$ws.Range("I113").Value = 4507.625
$ws.Range("J113").Value = 6729.706
$ws.Range("K113").Value = 4507.625
$ws.Range("L113").Value = 6729.706
$ws.Range("M113").Value = -1253.625
$ws.Range("N113").Value = -13237.706

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3445.6177
$ws.Range("I137").Value = 2303.1304
$ws.Range("K137").Value = 6909.3912
$ws.Range("M137").Value = -4359.3912

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 6406.136
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 6406.136
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 19218.408
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -29498.408

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 8181.6
$ws.Range("I141").Value = 9545.571
$ws.Range("K141").Value = 28636.713
$ws.Range("M141").Value = -23456.713

# ARM row 28
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 5311
$ws.Range("I28").Value = 2904.9285
$ws.Range("K28").Value = 2904.9285
$ws.Range("M28").Value = -2712.9285

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13441.517
$ws.Range("I32").Value = 10167.271
$ws.Range("J32").Value = 35542.668
$ws.Range("K32").Value = 10167.271
$ws.Range("L32").Value = 35542.668
$ws.Range("M32").Value = -9880.271000000001
$ws.Range("N32").Value = -36116.668

# ARM row 41
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 11351.333
$ws.Range("I41").Value = 11351.333
$ws.Range("K41").Value = 11351.333
$ws.Range("M41").Value = -10937.333

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3003.26
$ws.Range("J61").Value = 6321.8887
$ws.Range("L61").Value = 6321.8887
$ws.Range("N61").Value = -6745.8887

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 5579198.5
$ws.Range("J88").Value = 7597825.5
$ws.Range("L88").Value = 7597825.5
$ws.Range("N88").Value = -7598637.5

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 5579198.5
$ws.Range("J91").Value = 7597825.5
$ws.Range("L91").Value = 7597825.5
$ws.Range("N91").Value = -7600633.5

# ARM row 99
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H99").Value = 5311
$ws.Range("I99").Value = 2904.9285
$ws.Range("K99").Value = 2904.9285
$ws.Range("M99").Value = 90.07150000000001

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5760.2583
$ws.Range("I132").Value = 5571.273
$ws.Range("K132").Value = 16713.819
$ws.Range("M132").Value = -14183.819

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3003.26
$ws.Range("J136").Value = 6321.8887
$ws.Range("L136").Value = 18965.6661
$ws.Range("N136").Value = -24065.6661

# BSM row 81
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 36194.5
$ws.Range("J81").Value = 36194.5
$ws.Range("L81").Value = 36194.5
$ws.Range("N81").Value = -38316.5

# BSM row 84
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H84").Value = 36194.5
$ws.Range("J84").Value = 36194.5
$ws.Range("L84").Value = 108583.5
$ws.Range("N84").Value = -119191.5

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2936.2942
$ws.Range("I86").Value = 2780.5
$ws.Range("J86").Value = 3663.3333
$ws.Range("K86").Value = 2780.5
$ws.Range("L86").Value = 3663.3333
$ws.Range("M86").Value = -1657.5
$ws.Range("N86").Value = -5909.3333

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2936.2942
$ws.Range("I89").Value = 2780.5
$ws.Range("J89").Value = 3663.3333
$ws.Range("K89").Value = 13902.5
$ws.Range("L89").Value = 18316.6665
$ws.Range("M89").Value = -8286.5
$ws.Range("N89").Value = -29548.6665

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2775.4211
$ws.Range("I105").Value = 1715.8
$ws.Range("K105").Value = 1715.8
$ws.Range("M105").Value = 31.20000000000005

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 8871.875
$ws.Range("J62").Value = 8195
$ws.Range("L62").Value = 8195
$ws.Range("N62").Value = -9443

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 8871.875
$ws.Range("J65").Value = 8195
$ws.Range("L65").Value = 40975
$ws.Range("N65").Value = -47215

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3684.4092
$ws.Range("I132").Value = 2881.7058
$ws.Range("J132").Value = 6413.6
$ws.Range("K132").Value = 8645.117400000001
$ws.Range("L132").Value = 19240.8
$ws.Range("M132").Value = -6115.117400000001
$ws.Range("N132").Value = -24300.8

# CRP row 141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 207899
$ws.Range("J141").Value = 207899
$ws.Range("L141").Value = 207899
$ws.Range("N141").Value = -218259

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 8584.615
$ws.Range("I5").Value = 614.6
$ws.Range("K5").Value = 1843.8
$ws.Range("M5").Value = -1731.8

# CUL row 116
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 3574.25
$ws.Range("J116").Value = 3574.25
$ws.Range("L116").Value = 10722.75
$ws.Range("N116").Value = -17606.75

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3596.7144
$ws.Range("I131").Value = 2694.4614
$ws.Range("J131").Value = 4378.6665
$ws.Range("K131").Value = 8083.3842
$ws.Range("L131").Value = 13135.9995
$ws.Range("M131").Value = -3043.3842
$ws.Range("N131").Value = -23215.9995

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3375.7144
$ws.Range("J132").Value = 5999.6665
$ws.Range("L132").Value = 53996.9985
$ws.Range("N132").Value = -59056.9985

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 8584.615
$ws.Range("I135").Value = 614.6
$ws.Range("K135").Value = 5531.400000000001
$ws.Range("M135").Value = -2996.400000000001

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 268.82352
$ws.Range("I107").Value = 300.76923
$ws.Range("J107").Value = 165
$ws.Range("K107").Value = 300.76923
$ws.Range("L107").Value = 165
$ws.Range("M107").Value = 1619.23077
$ws.Range("N107").Value = -4005

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9041.362999999999
$ws.Range("I40").Value = 9553.333000000001
$ws.Range("J40").Value = 7944.2856
$ws.Range("K40").Value = 9553.333000000001
$ws.Range("L40").Value = 7944.2856
$ws.Range("M40").Value = -9417.333000000001
$ws.Range("N40").Value = -8216.285599999999

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4831.3022
$ws.Range("I132").Value = 4484.25
$ws.Range("J132").Value = 5269.684
$ws.Range("K132").Value = 13452.75
$ws.Range("L132").Value = 15809.052
$ws.Range("M132").Value = -10922.75
$ws.Range("N132").Value = -20869.052

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5202.483
$ws.Range("I136").Value = 3268.8235
$ws.Range("J136").Value = 7941.8335
$ws.Range("K136").Value = 9806.470499999999
$ws.Range("L136").Value = 23825.5005
$ws.Range("M136").Value = -7256.470499999999
$ws.Range("N136").Value = -28925.5005

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 27779560
$ws.Range("I62").Value = 3125
$ws.Range("J62").Value = 37038372
$ws.Range("K62").Value = 3125
$ws.Range("L62").Value = 37038372
$ws.Range("M62").Value = -2501
$ws.Range("N62").Value = -37039620

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 27779560
$ws.Range("I65").Value = 3125
$ws.Range("J65").Value = 37038372
$ws.Range("K65").Value = 15625
$ws.Range("L65").Value = 185191860
$ws.Range("M65").Value = -12505
$ws.Range("N65").Value = -185198100

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4615.7
$ws.Range("I122").Value = 4615.7
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 13847.1
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -11397.1
$ws.Range("N122").ClearContents()

# WVR row 135
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 60384.54
$ws.Range("J135").Value = 60384.54
$ws.Range("L135").Value = 60384.54
$ws.Range("N135").Value = -70524.54000000001
